$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the floating point rounding issue: 0.6+0.3 (=0.89999999999999991) -> 0.1+0.7 (=0.79999999999999993)
$ws.Range("A2").Formula = "=0.1+0.7"

# Update the active selection to A2 (was A3)
$ws.Range("A2").Select()
